# Generate Report for Archive
#
# The localization status report is regenerated: the entry for
# "19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md" (previously "Ready for handoff")
# moves down past "6905292b-5e95-4e1f-a663-afb2b2ba38e1.md" and
# "bdb88f2c-b595-470e-976f-47c3662a2aed.md", which both flip to
# "In Translation" status (they now precede it in the report).
#
# This touches the "Overview" sheet (rows 3-5, cols A/B/E/F) and the
# "zh-cn" / "de-de" sheets (rows 3-5, cols A/C/G), plus the hyperlink
# display text that tracks the file name shown in column A/B of each row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = "6905292b-5e95-4e1f-a663-afb2b2ba38e1.md"
$ws.Range("B3").Value = "e2e\6905292b-5e95-4e1f-a663-afb2b2ba38e1.md"
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"

$ws.Range("A4").Value = "bdb88f2c-b595-470e-976f-47c3662a2aed.md"
$ws.Range("B4").Value = "e2e\bdb88f2c-b595-470e-976f-47c3662a2aed.md"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"

$ws.Range("A5").Value = "19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md"
$ws.Range("B5").Value = "e2e\19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md"

# Rebuild the hyperlinks so rId/order + display text matches the new
# row contents (the link target for a given row position is unchanged,
# only the visible display text tracks the file now shown there).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c3f6bb9311a84d87da127f67b2b53a9f1051d66/e2e/08251479-1d8e-4084-a7af-8626ee35bf73.md", [Type]::Missing, [Type]::Missing, "e2e\08251479-1d8e-4084-a7af-8626ee35bf73.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dace41f88439ff83eeba76491264b39f7ba16ad5/e2e/19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md", [Type]::Missing, [Type]::Missing, "e2e\6905292b-5e95-4e1f-a663-afb2b2ba38e1.md")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e530747320c221415b1bfa2b4ddb51cd6abfa028/e2e/6905292b-5e95-4e1f-a663-afb2b2ba38e1.md", [Type]::Missing, [Type]::Missing, "e2e\bdb88f2c-b595-470e-976f-47c3662a2aed.md")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e530747320c221415b1bfa2b4ddb51cd6abfa028/e2e/bdb88f2c-b595-470e-976f-47c3662a2aed.md", [Type]::Missing, [Type]::Missing, "e2e\19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99387b10abdf407445a92dd2ca9f6e445eb06879/e2e/cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md", [Type]::Missing, [Type]::Missing, "e2e\cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f547dbfa1528a922b84290cb8320a2c1c9f66dc/e2e/d86e2b24-8ab0-4acd-9e6c-e3e1822202f5.md", [Type]::Missing, [Type]::Missing, "e2e\d86e2b24-8ab0-4acd-9e6c-e3e1822202f5.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = "6905292b-5e95-4e1f-a663-afb2b2ba38e1.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("G3").Value = "6905292b-5e95-4e1f-a663-afb2b2ba38e1.612a408cd7b71953c6b893a10de8a67d43ed8f6e.zh-cn.xlf"

$ws.Range("A4").Value = "bdb88f2c-b595-470e-976f-47c3662a2aed.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "bdb88f2c-b595-470e-976f-47c3662a2aed.1dc97912d6dd67b06405175a1e01c8b7eebab771.zh-cn.xlf"

$ws.Range("A5").Value = "19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md"
$ws.Range("G5").Value = "19ed8b35-cfdf-49a5-8db7-4708c4e08b75.aee16b197219b2ff05a642e563810e770a1f873f.zh-cn.xlf"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c3f6bb9311a84d87da127f67b2b53a9f1051d66/e2e/08251479-1d8e-4084-a7af-8626ee35bf73.md", [Type]::Missing, [Type]::Missing, "08251479-1d8e-4084-a7af-8626ee35bf73.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6ee572d90477dbd0771a1ec6350988935566f0ff/e2e/08251479-1d8e-4084-a7af-8626ee35bf73.md", [Type]::Missing, [Type]::Missing, "08251479-1d8e-4084-a7af-8626ee35bf73.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dace41f88439ff83eeba76491264b39f7ba16ad5/e2e/19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md", [Type]::Missing, [Type]::Missing, "6905292b-5e95-4e1f-a663-afb2b2ba38e1.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e530747320c221415b1bfa2b4ddb51cd6abfa028/e2e/6905292b-5e95-4e1f-a663-afb2b2ba38e1.md", [Type]::Missing, [Type]::Missing, "bdb88f2c-b595-470e-976f-47c3662a2aed.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e530747320c221415b1bfa2b4ddb51cd6abfa028/e2e/bdb88f2c-b595-470e-976f-47c3662a2aed.md", [Type]::Missing, [Type]::Missing, "19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99387b10abdf407445a92dd2ca9f6e445eb06879/e2e/cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md", [Type]::Missing, [Type]::Missing, "cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md")
$ws.Hyperlinks.Add($ws.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6ee572d90477dbd0771a1ec6350988935566f0ff/e2e/cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md", [Type]::Missing, [Type]::Missing, "cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f547dbfa1528a922b84290cb8320a2c1c9f66dc/e2e/d86e2b24-8ab0-4acd-9e6c-e3e1822202f5.md", [Type]::Missing, [Type]::Missing, "d86e2b24-8ab0-4acd-9e6c-e3e1822202f5.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = "6905292b-5e95-4e1f-a663-afb2b2ba38e1.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("G3").Value = "6905292b-5e95-4e1f-a663-afb2b2ba38e1.612a408cd7b71953c6b893a10de8a67d43ed8f6e.de-de.xlf"

$ws.Range("A4").Value = "bdb88f2c-b595-470e-976f-47c3662a2aed.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "bdb88f2c-b595-470e-976f-47c3662a2aed.1dc97912d6dd67b06405175a1e01c8b7eebab771.de-de.xlf"

$ws.Range("A5").Value = "19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md"
$ws.Range("G5").Value = "19ed8b35-cfdf-49a5-8db7-4708c4e08b75.aee16b197219b2ff05a642e563810e770a1f873f.de-de.xlf"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c3f6bb9311a84d87da127f67b2b53a9f1051d66/e2e/08251479-1d8e-4084-a7af-8626ee35bf73.md", [Type]::Missing, [Type]::Missing, "08251479-1d8e-4084-a7af-8626ee35bf73.md")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e6f94f9f0648615fd8f8899d809f08ddf99c0d5c/e2e/08251479-1d8e-4084-a7af-8626ee35bf73.md", [Type]::Missing, [Type]::Missing, "08251479-1d8e-4084-a7af-8626ee35bf73.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dace41f88439ff83eeba76491264b39f7ba16ad5/e2e/19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md", [Type]::Missing, [Type]::Missing, "6905292b-5e95-4e1f-a663-afb2b2ba38e1.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e530747320c221415b1bfa2b4ddb51cd6abfa028/e2e/6905292b-5e95-4e1f-a663-afb2b2ba38e1.md", [Type]::Missing, [Type]::Missing, "bdb88f2c-b595-470e-976f-47c3662a2aed.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e530747320c221415b1bfa2b4ddb51cd6abfa028/e2e/bdb88f2c-b595-470e-976f-47c3662a2aed.md", [Type]::Missing, [Type]::Missing, "19ed8b35-cfdf-49a5-8db7-4708c4e08b75.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/99387b10abdf407445a92dd2ca9f6e445eb06879/e2e/cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md", [Type]::Missing, [Type]::Missing, "cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md")
$ws.Hyperlinks.Add($ws.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e6f94f9f0648615fd8f8899d809f08ddf99c0d5c/e2e/cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md", [Type]::Missing, [Type]::Missing, "cc5fb36c-4deb-4ff7-9c07-8ce5503f3a78.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f547dbfa1528a922b84290cb8320a2c1c9f66dc/e2e/d86e2b24-8ab0-4acd-9e6c-e3e1822202f5.md", [Type]::Missing, [Type]::Missing, "d86e2b24-8ab0-4acd-9e6c-e3e1822202f5.md")
